$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = '@'
$cell.Value = '43.748.18'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.246.03'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.34%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.06%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 2)
$cell.NumberFormat = '@'
$cell.Value = 'XRP'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.645'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.94%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 2)
$cell.NumberFormat = '@'
$cell.Value = 'BNB'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 3)
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '230.29'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.34%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '64.29'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +5.11%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.02%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.448'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +5.84%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0976'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.62%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '56.92'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.84%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '26.67'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +12.54%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.105'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.27%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.581.92'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.12%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.07'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.39%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.830'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.77%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.254.38'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.86%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '43.665.43'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.29%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0985'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +5.69%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '73.11'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.01%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.02'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.03%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '250.23'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.13%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -5.03%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.33%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.33'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +24.23%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '9.97'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.10%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '170.59'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.17%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.58%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '20.79'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +1.39%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.00%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.07%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0701'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +6.63%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.76'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.30%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.86'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.86%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.70'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +2.83%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.45'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.72%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.28'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -4.52%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.21%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.11%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.000224'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.35%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '17.27'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +3.47%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0962'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.14%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '8.16'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -6.56%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.20'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -0.89%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '97.21'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -1.20%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -2.33%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.36'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +6.10%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '10.00'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = '@'
$cell.Value = '  +4.67%  '
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.432.40'
$cell.Style = 'Normal'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = '@'
$cell.Value = '  -3.09%  '
$cell.Style = 'Normal'
